$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.593.98'
$ws.Range("E2").Value = '  -3.36%  '
$ws.Range("D3").Value = '1.850.68'
$ws.Range("E3").Value = '  -4.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -1.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '336.59'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4657'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.44%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3911'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.70%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.12'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07908'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9828'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.34'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -6.45%  '
$ws.Range("D13").Value = '1.870.70'
$ws.Range("E13").Value = '  -3.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.851'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -4.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.016'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06906'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.59'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001004'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.13'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.87%  '
$ws.Range("D22").Value = '28.605.45'
$ws.Range("E22").Value = '  -3.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.408'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.34'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -5.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.143'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.60%  '
$ws.Range("D26").Value = '2.091.42'
$ws.Range("E26").Value = '  -3.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.67'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.44'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.079'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -5.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.028'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.69'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9733'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09393'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.367'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.484'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.349'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06163'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02200'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.162'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5723'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.705'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.82%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.19'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -5.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1801'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.356'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -4.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.254'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5397'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.73'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -5.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.07141'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -4.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.910'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '115.27'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.06'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.17%  '
